$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 7 (apportionment / states reprocessing) with new values
$ws.Range("B7").Value = 0.09255869550917829
$ws.Range("C7").Value = 0.07569922773292695
$ws.Range("D7").Value = 0.08300328292707819
$ws.Range("E7").Value = 0.04509469448561292
